$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A and append two new tracker rows
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$row1 = $lastRow + 1
$row2 = $lastRow + 2

$ws.Cells.Item($row1, 1).Value = "G1"
$ws.Cells.Item($row1, 2).Value = "Test1"
$ws.Cells.Item($row1, 3).Value = 45895
$ws.Cells.Item($row1, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($row1, 4).Value = 0.7201030745259143
$ws.Cells.Item($row1, 5).Value = 0
$ws.Cells.Item($row1, 6).Value = -0.01

$ws.Cells.Item($row2, 1).Value = "G2"
$ws.Cells.Item($row2, 2).Value = "sedrftgyhuioygtfrd"
$ws.Cells.Item($row2, 3).Value = 45895
$ws.Cells.Item($row2, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($row2, 4).Value = 0.7201030745259143
$ws.Cells.Item($row2, 5).Value = 0
$ws.Cells.Item($row2, 6).Value = -0.01
